$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-02-18T21:06:08.465605"

# Row 2
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0.00006944444444444444
$ws.Range("K2").Value = 476
$ws.Range("L2").Value = 0.000952
$ws.Range("M2").Value = $newTimestamp

# Row 3
$ws.Range("K3").Value = 472
$ws.Range("L3").Value = 0.00472
$ws.Range("M3").Value = $newTimestamp

# Row 4
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = $newTimestamp

# Row 5
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = $newTimestamp

# Row 6 (only timestamp changes)
$ws.Range("M6").Value = $newTimestamp

# Row 7 (only timestamp changes)
$ws.Range("M7").Value = $newTimestamp

# Row 8 (only timestamp changes)
$ws.Range("M8").Value = $newTimestamp

# Row 9
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = $newTimestamp

# Row 10
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = $newTimestamp

# Row 11
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = $newTimestamp

# Row 12
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = $newTimestamp

# Row 13 (only timestamp changes)
$ws.Range("M13").Value = $newTimestamp

# Row 14
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = $newTimestamp

$wb.Save()
